$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (row 38) was appended, duplicating the last existing
# data row (row 37) -- same date/price text and same cell styles.
$ws.Range("A37:B37").Copy() | Out-Null
$ws.Range("A38:B38").PasteSpecial(-4104) | Out-Null  # xlPasteAll
$excel.CutCopyMode = $false
